$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mattress pad")
$arr = New-Object 'object[,]' 55,2
$arr[0,0] = 'Product Name'
$arr[0,1] = 'Star Rank'
$arr[1,0] = '[Sponsored]Mattress Pad Cover Queen Size Pillowtop 300TC Down Alternative Mattress Topper with 8-21-Inch Deep Pocket'
$arr[1,1] = '1,1,1'
$arr[2,0] = '[Sponsored]Queen Overfilled Mattress Pad Cover 8-21”Deep Pocket-Cooling Mattress Topper Snow Down Alternative'
$arr[2,1] = '1,1,2'
$arr[3,0] = '[Sponsored]Fairyland Mattress Pad Cover with 300TC 100% Cotton Quilted Down Alternative Filled Mattress Topper,8-21 Inch Deep Pocket (TwinXL)'
$arr[3,1] = '1,1,3'
$arr[4,0] = 'Quilted Fitted Mattress Pad (Queen) - Mattress Cover Stretches up to 16 Inches Deep - Mattress Topper by Utopia Bedding'
$arr[4,1] = '1,2,1'
$arr[5,0] = 'Zen Bamboo Ultra Soft Fitted Bamboo Mattress Pad - Premium Hypoallergenic Bamboo Mattress Topper with Honeycomb Cooling Technology - Queen'
$arr[5,1] = '1,2,2'
$arr[6,0] = 'Beckham Hotel Collection Premium Microplush Mattress Pad - Hypoallergenic Ultra Soft Overfilled Topper with Deep Fit - Queen'
$arr[6,1] = '1,2,3'
$arr[7,0] = 'Mattress Pad Cover-Cotton Top with Stretches to 18” Deep Pocket Fits Up to 8”-21” Cooling White Bed Topper (Down Alternative, Queen)'
$arr[7,1] = '1,3,1'
$arr[8,0] = 'Queen Overfilled Mattress Pad Cover 8-21”Deep Pocket-Cooling Mattress Topper Snow Down Alternative'
$arr[8,1] = '1,3,2'
$arr[9,0] = 'Mattress Pad Cover Queen Size Pillowtop 300TC Down Alternative Mattress Topper with 8-21-Inch Deep Pocket'
$arr[9,1] = '1,3,3'
$arr[10,0] = 'Beckham Hotel Collection Microfiber Mattress Pad - Quilted, Hypoallergnic, and Water-Resistant - Queen'
$arr[10,1] = '1,4,1'
$arr[11,0] = 'Hypoallergenic Quilted Stretch-to-Fit Mattress Pad By Hanna Kay, 10 Year Warranty-Clyne Collection (Queen)'
$arr[11,1] = '1,4,2'
$arr[12,0] = 'Queen Mattress Pad Cover 8-21”Deep Pocket - Cooling Mattress Topper Overfilled 300TC Snow Down Alternative'
$arr[12,1] = '1,4,3'
$arr[13,0] = 'Premium Mattress Pad (Queen) - Quilted Fitted Mattress Topper Stretches Upto 15 Inches Deep - Plush and Soft Mattress Protector And Cover With Deep Pockets By Utopia Bedding'
$arr[13,1] = '1,5,1'
$arr[14,0] = 'INGALIK Hotel Luxury Collection Quilted Fitted Mattress Topper Down Alternative Overfilled Mattress Pad Bed Cover Stretches up to 21 Inches Deep by (Queen 60x80x18inch)'
$arr[14,1] = '1,5,2'
$arr[15,0] = 'Beckham Hotel Collection 1500 Series Microfiber Mattress Pad - Quilted, Hypoallergnic, and Water-Resistant - King'
$arr[15,1] = '1,5,3'
$arr[16,0] = 'Mattress Pad Cover with 18" Deep Pocket 300TC Cotton Down Mattress Topper By HYPNOS Mattress Topper Hypoallergenic Quilted Stretch-to-Fit,King'
$arr[16,1] = '1,6,1'
$arr[17,0] = 'Maevis Mattress Pad Cover 100% 300TC Cotton with 8-21 Inch Deep Pocket White Overfilled Bed Mattress Topper (Down Alternative, Twin XL)'
$arr[17,1] = '1,6,2'
$arr[18,0] = 'Sleep Restoration Fitted Microfiber Mattress Pad Cover - Plush Quilted Luxurious Mattress Topper - Queen'
$arr[18,1] = '1,6,3'
$arr[19,0] = 'Mattress Pad Cover (King Size)- Cooling Mattress Topper with Thick Cotton 8-21-Inch Deep Pocket - Quilted Fitted Pillowtop by Sonoro Kate'
$arr[19,1] = '1,7,1'
$arr[20,0] = 'Mattress Pad Cover with 18” Deep Pocket Overfilled 100% 300TC Cotton White Bed Topper By WarmHarbor Mattress Topper (Down Alternative, Full)'
$arr[20,1] = '1,7,2'
$arr[21,0] = 'Pillowtop Mattress Pad with Fitted Skirt - Extra Plush Topper Found in Marriott Hotels - Made in the USA, Queen'
$arr[21,1] = '1,7,3'
$arr[22,0] = 'Fairyland Mattress Pad Cover with 300TC 100% Cotton Quilted Down Alternative Filled Mattress Topper,8-21 Inch Deep Pocket (TwinXL)'
$arr[22,1] = '1,8,1'
$arr[23,0] = 'Pillowtop Mattress Pad Cover Queen Size - Hypoallergenic - Cotton Down Alternative Filled Mattress Topper'
$arr[23,1] = '1,8,2'
$arr[24,0] = 'Bamboo Overfilled Pillow Top Mattress Pad | Superb Temperature Regulation | Made in the USA, King'
$arr[24,1] = '1,8,3'
$arr[25,0] = 'Mattress Pad Full Size with 24 inch Deep Pocket Microplush Mattress Topper with Fitted Skirt Quilted Stretch Pillow Top by Naluka（54”x75”)'
$arr[25,1] = '1,9,1'
$arr[26,0] = 'Luxurious - Fitted Down Alternative Mattress Pad - 100% Cotton Top Mattress Topper, 300 Thread Count - Mattress Cover Stretches Up to 16 - Full Size (54x75")'
$arr[26,1] = '1,9,2'
$arr[27,0] = 'Twin XL Mattress Pad - 300TC Down Alternative Pillow Top Mattress Topper,Quilted Deep Pocket Fitted Mattress Cover (8"-21")'
$arr[27,1] = '1,9,3'
$arr[28,0] = 'Mattress Pad King Size 400TC Cotton Top 3M Water Resistant Hypoallergenic-71oz Down Alternative Filling Pillowtop Mattress Topper Cover-Fitted Quilted 8-21 Inch Deep Pocket'
$arr[28,1] = '1,10,1'
$arr[29,0] = 'Mattress Pad Cover, Microfiber, Soft, Hypoallergenic, Mattress Topper with Deep Pocket(Queen,Superior)'
$arr[29,1] = '1,10,2'
$arr[30,0] = 'Rayon from Bamboo Extra Thick Mattress Pad with Fitted Skirt - Extra Plush Cooling Topper - Hypoallergenic - Proudly Made in the USA, Twin XL'
$arr[30,1] = '1,10,3'
$arr[31,0] = 'Amazon recommendation'
$arr[31,1] = '1,11,1'
$arr[32,0] = 'Mattress Pad Full Size Hypoallergenic - Antibacterial, Breathable - Ultra Soft Quilted Mattress Protector, Fitted Sheet Mattress Cover White by Bedsure'
$arr[32,1] = '1,11,2'
$arr[33,0] = 'Quilted Fitted Mattress Pad (King)-Mattress Cover Stretches up 8-21" Deep Pocket Down Alternative Filling Mattress Topper'
$arr[33,1] = '1,11,3'
$arr[34,0] = 'Naturepedic Organic Waterproof Fitted Stretch Knit Protector Pad - Twin'
$arr[34,1] = '1,12,1'
$arr[35,0] = 'Pillow-Top Premium Mattress Pad - 1.5 Inch Cooling Down Alternative Polygel Filled Microplush Super-Soft Hypoallergenic Topper (Twin XL/Twin Extra Long)'
$arr[35,1] = '1,12,2'
$arr[36,0] = 'SleepJoy 3" ViscO2 Ventilated Memory Foam Mattress Topper, Queen'
$arr[36,1] = '1,12,3'
$arr[37,0] = 'Sleep Philosophy All Natural Cotton Filled Mattress Pad, Queen'
$arr[37,1] = '1,13,1'
$arr[38,0] = 'AmazonBasics Hypoallergenic Quilted Mattress Topper, 18" Deep, Full'
$arr[38,1] = '1,13,2'
$arr[39,0] = 'The Grand Fitted Quilted Mattress Pad Cover Hypoallergenic (Stretches to 18" Deep - Queen - 60x80") Queen Size Mattress Protector'
$arr[39,1] = '1,13,3'
$arr[40,0] = 'Queen Fitted Quilted Mattress Pad Cover 8-21”Deep Pocket-Down Alternative Mattress Topper'
$arr[40,1] = '1,14,1'
$arr[41,0] = 'Micropuff Down Alternative Mattress Pad - White Quilted Fitted Mattress Topper (Twin Size - 39"x75") Microfiber Mattress Cover Stretches up to 15"'
$arr[41,1] = '1,14,2'
$arr[42,0] = 'Bamboo Mattress Pad with Fitted Skirt - Extra Plush Cooling Topper - Hypoallergenic - Made in the USA, Full'
$arr[42,1] = '1,14,3'
$arr[43,0] = 'Hospitology Heavenly Microfiber Goose Down Alternative Overstuffed Hypoallergenic Mattress Pad / Topper, 54-Inch by 75-Inch, Full/Double'
$arr[43,1] = '1,15,1'
$arr[44,0] = 'Shilucheng Mattress Pad Queen Size Ultra Soft Rayon Derived from Plush - 8-21-Inch Deep Pocket Premium Hypoallergenic Mattress Topper'
$arr[44,1] = '1,15,2'
$arr[45,0] = 'Happsy Organic Mattress Protector Pad - Twin'
$arr[45,1] = '1,15,3'
$arr[46,0] = 'Merous Queen Size Cotton Mattress Pad Down Alternative Mattress Cover - Hypoallergenic Fitted Quilted Mattress Topper - Stretches up to 18 Inches Deep'
$arr[46,1] = '1,16,1'
$arr[47,0] = 'Superior Queen Size Premium 100% Waterproof Mattress Protector Pad - 100% Cotton Terry Surface, Hypoallergenic, Deep Pocket Skirt Fits Up to 22" Mattress, 15-Year Warranty'
$arr[47,1] = '1,16,2'
$arr[48,0] = 'Red Nomad - Queen Size 2 Inch Thick, Ultra Premium Visco Elastic Memory Foam Mattress Pad Bed Topper - Made in the USA'
$arr[48,1] = '1,16,3'
$arr[49,0] = 'RV Mattress Pad - Extra Plush Topper with Fitted Skirt - Found in Marriott Hotels - Made in the USA - Hypoallergenic - Mattress Cover for RV, Camper - Short Queen'
$arr[49,1] = '1,17,1'
$arr[50,0] = 'Allrange 300TC Cool Tencel Hypoallergenic Quilted Mattress Pad, Stretch-up-to 22", Fitted Tencel Polyester Fill, Silky Cotton Tencel Cover,OEKO-TEX Certified, Queen'
$arr[50,1] = '1,17,2'
$arr[51,0] = 'Pressure Relief Mattress Pad with Fitted Skirt |Bedsore Prevention Mattress Pads | Hypoallergenic Mattress Topper | Made in the USA, Queen'
$arr[51,1] = '1,17,3'
$arr[52,0] = '[Sponsored]Queen Mattress Pad Cover 8-21”Deep Pocket - Cooling Mattress Topper Overfilled 300TC Snow Down Alternative'
$arr[52,1] = '1,18,1'
$arr[53,0] = '[Sponsored]Bamboo Overfilled and ExtraThick 1-Piece Pillow Top Mattress Pad, Twin XL'
$arr[53,1] = '1,18,2'
$arr[54,0] = '[Sponsored]Pillowtop Mattress Pad with Fitted Skirt - Extra Plush Topper Found in Marriott Hotels - Made in the USA, Twin XL'
$arr[54,1] = '1,18,3'
$ws.Range("A1:B55").Value = $arr
Write-Host "Done writing mattress pad rows"
